$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1 currently repeats the first data row by mistake; turn it into a
#     proper header row (same column names used on the 土地/建物 sheets),
#     and extend it with the property_category..index columns. ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Apply the existing header formatting (bold/border/centered) to the
# whole header row, including the newly-added H1:N1 cells.
$ws.Range("B1").Copy()
$ws.Range("B1:N1").PasteSpecial(-4122)

# --- Row 3: the capacity value for the MAZDAMPV row was stored as the
#     text "3000" - normalize it to a real number like the other rows. ---
$ws.Cells.Item(3, 3).Value = 3000

# --- Rows 2-4: append property_category / category / date /
#     legislator_name / legislator_id / source_file / index, mirroring the
#     columns already present on the 土地 / 建物 sheets. "date" must stay
#     text (not get auto-parsed into a serial date number). ---
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 8).Value = "land"
    $ws.Cells.Item($r, 9).Value = "normal"

    $ws.Cells.Item($r, 10).NumberFormat = "@"
    $ws.Cells.Item($r, 10).Value = "2013-12-02"

    $ws.Cells.Item($r, 11).Value = "李俊俋"
    $ws.Cells.Item($r, 12).Value = 1738
    $ws.Cells.Item($r, 13).Value = "tmp52b51"
    $ws.Cells.Item($r, 14).Value = 28 + $r

    # Re-apply the plain data-row formatting (matches columns B:G on this
    # row) now that the date column's number format was changed above.
    $src = "G" + $r
    $dst = "H" + $r + ":N" + $r
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial(-4122)
}
